$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values ("310.41", "8.46", ...)
# as literal text in the source workbook, not as numbers. Plain assignment via
# .Value would let Excel auto-convert those strings to real numbers, so for the
# cells whose new value would parse as a number we force a text number format
# first (per-cell, since NumberFormat on a multi-area Range only hit the first
# area here), write the value, then restore General/Normal so no cell is left
# with stray formatting.
$textForceCells = @("D5","D6","D10","D12","D14","D16","D19","D22","D23","D24","D25","D27","D28","D29","D30","D31","D33","D34","D36","D37","D38","D39","D40","D42","D44","D45","D47","D48","D49","D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.155.49"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "2.359.18"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "310.41"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "108.41"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").Value = "41.29"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "8.46"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "0.981"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "2.719.39"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "15.23"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "2.356.58"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "45.165.96"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").Value = "14.19"
$ws.Range("E19").Value = "  +7.96%  "
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "73.11"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "259.87"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "11.15"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").Value = "  -5.02%  "
$ws.Range("D29").Value = "2.34"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("D30").Value = "0.0963"
$ws.Range("E30").Value = "  +9.29%  "
$ws.Range("D31").Value = "22.31"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "168.71"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  +7.40%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("D37").Value = "4.83"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").Value = "2.96"
$ws.Range("E38").Value = "  +4.44%  "
$ws.Range("D39").Value = "3.92"
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").Value = "0.0355"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("D42").Value = "99.42"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "69.48"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "12.87"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "81.39"
$ws.Range("E47").Value = "  +4.68%  "
$ws.Range("D48").Value = "112.40"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "5.52"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "9.26"
$ws.Range("E50").Value = "  +4.66%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.667.25"
$ws.Range("E51").Value = "  +0.56%  "

# Restore default number format/style on the cells we force-formatted above.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
